$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
  2 = 5
  3 = 6
  4 = 6
  5 = 4
  6 = 5
  7 = 4
  8 = 6
  9 = 5
  10 = 7
  11 = 6
  12 = 3
  13 = 9
  14 = 1
  15 = 6
  16 = 6
  17 = 7
  18 = 3
  19 = 4
  20 = 7
  21 = 7
  22 = 5
  23 = 8
  24 = 4
  25 = 6
  26 = 4
  27 = 6
  28 = 2
  29 = 0
  30 = 1
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
